# "correcciones generales de contenido estadisticas y variables"
#
# This script reproduces, via PowerPoint COM automation, the small
# layout/content corrections recorded in the target diff:
#   - the three slide-master footer / slide-number / date placeholders
#     are nudged to a very slightly smaller size;
#   - the background picture and the certificate text boxes on slide 1
#     are repositioned / resized;
#   - the "Rectangulo 6" text box (company-name placeholder) loses its
#     leading space run and is moved/resized to sit next to the label.
#
# NOTE on precision: PowerPoint stores Shape.Left/Top/Width/Height in
# points as 32-bit floats, while the OOXML stores EMU (1 pt = 12700 EMU)
# as integers. round-tripping pt -> EMU truncates, so literal
# `target_emu / 12700.0` can miss the exact target EMU by 1 unit. The
# point literals below were solved so that, after the float32 round
# trip performed by the host, they serialize back to precisely the EMU
# values required by the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sm = $p.SlideMaster

# ---------------------------------------------------------------------
# Slide master: footer / slide-number / date placeholders shrink by a
# hair (cx -720 EMU, cy -720 EMU) with their position unchanged.
# ---------------------------------------------------------------------

$ftr = $sm.Shapes.Item(1)        # PlaceHolder 1 - ftr
$ftr.Width = 242.84409448818897
$ftr.Height = 28.57322834645669

$sldNum = $sm.Shapes.Item(2)     # PlaceHolder 2 - sldNum
$sldNum.Width = 161.8299212598425
$sldNum.Height = 28.57322834645669

$dt = $sm.Shapes.Item(3)         # PlaceHolder 3 - dt
$dt.Width = 161.8299212598425
$dt.Height = 28.57322834645669

# ---------------------------------------------------------------------
# Slide 1 shapes
# ---------------------------------------------------------------------

# Imagen 2 (background picture) - position unchanged, tiny resize
$img = $s.Shapes.Item(1)
$img.Width = 719.8299255598425
$img.Height = 539.8299255598425

# Rectangulo 6 - "(Nombre_Comercio)" label: drop the leading space run,
# then move/resize the box to its new place next to "Comercio Afiliado:"
$r6 = $s.Shapes.Item(2)
$r6.TextFrame.TextRange.Characters(1, 1).Text = ""
$r6.Left = 252.62362204724408
$r6.Top = 167.896064792126
$r6.Width = 197.74488188976378
$r6.Height = 31.06771653543307

# Rectangulo 8 - "RIF:" label, moves slightly
$r8 = $s.Shapes.Item(3)
$r8.Left = 273.5716535433071
$r8.Top = 200.1543351086614
$r8.Width = 51.70393700787402
$r8.Height = 35.85826771653543

# Rectangulo 9, moves slightly
$r9 = $s.Shapes.Item(4)
$r9.Left = 92.46614463228346
$r9.Top = 266.25826771653544
$r9.Width = 95.04566929133858
$r9.Height = 28.658267716535434

# CuadroTexto 3 - position unchanged, tiny resize
$ct3 = $s.Shapes.Item(5)
$ct3.Width = 98.2771653543307
$ct3.Height = 28.658267716535434

# CuadroTexto 4 - position unchanged, tiny resize
$ct4 = $s.Shapes.Item(6)
$ct4.Width = 97.90866471732284
$ct4.Height = 21.42992125984252

# Rectangulo 5 - position unchanged, tiny resize
$r5 = $s.Shapes.Item(7)
$r5.Width = 127.36063012125985
$r5.Height = 105.98740157480314

# CuadroTexto 1 - position unchanged, tiny resize
$ct1 = $s.Shapes.Item(8)
$ct1.Width = 384.57637795275593
$ct1.Height = 28.658267716535434

# CuadroTexto 7 - position unchanged, tiny resize
$ct7 = $s.Shapes.Item(9)
$ct7.Width = 139.66299212598426
$ct7.Height = 28.658267716535434

Write-Output "Applied general content/layout corrections."
